$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'android.view.View; android.os.Bundle; android.widget.ImageView; android.widget.ListView; android.view.ViewGroup; com.cjj.MaterialRefreshListener; android.support.v7.widget.LinearLayoutManager; android.view.MenuItem; android.graphics.Color; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; android.support.v7.widget.Toolbar; java.util.Random; android.view.Menu; android.support.v7.widget.DefaultItemAnimator; android.support.v7.widget.RecyclerView; android.widget.Toast; java.util.List; java.util.ArrayList; android.content.Context; android.content.Intent; android.support.annotation.Nullable'
$ws.Range("B3").Value = 'android.animation.ValueAnimator; android.view.View; android.view.animation.LinearInterpolator; android.animation.ObjectAnimator; android.graphics.Color; android.util.AttributeSet; android.support.v4.view.ViewCompat; android.widget.FrameLayout; android.content.Context'
$ws.Range("B4").Value = 'android.view.View; android.view.animation.AccelerateDecelerateInterpolator; android.graphics.RadialGradient; android.graphics.Paint.Style; android.widget.AbsListView; com.cjj.header.SunLayout; android.os.Bundle; android.graphics.PixelFormat; android.widget.ImageView; android.widget.ListView; android.util.TypedValue; android.graphics.drawable.ShapeDrawable; com.cjj.MaterialRefreshListener; android.view.WindowManager; android.view.ViewGroup; android.support.v7.widget.LinearLayoutManager; android.graphics.Rect; android.os.Build; android.os.Handler; android.view.MenuItem; android.support.v4.view.ViewPropertyAnimatorUpdateListener; android.graphics.Color; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; android.support.v4.view.ViewPropertyAnimatorCompat; android.view.animation.BounceInterpolator; android.support.v7.widget.Toolbar; android.util.DisplayMetrics; java.util.Random; com.cjj.MaterialHeadListener; android.view.animation.DecelerateInterpolator; android.animation.ValueAnimator; android.view.animation.Transformation; android.view.Menu; android.util.Log; android.support.v7.widget.DefaultItemAnimator; android.view.animation.LinearInterpolator; android.graphics.ColorFilter; android.content.res.Resources; android.support.v7.widget.RecyclerView; android.net.Uri; android.widget.ArrayAdapter; android.view.animation.Interpolator; android.widget.Toast; java.util.List; java.util.ArrayList; android.content.Context; android.graphics.PaintFlagsDrawFilter; android.support.v7.app.AppCompatActivity; android.content.Intent; android.graphics.Shader; android.graphics.drawable.Drawable; android.graphics.drawable.Animatable; android.graphics.drawable.shapes.OvalShape; android.widget.LinearLayout; android.animation.ObjectAnimator; android.util.AttributeSet; android.content.res.TypedArray; android.view.Gravity; android.view.MotionEvent; android.support.v4.view.ViewCompat; android.support.annotation.Nullable; android.widget.FrameLayout; android.view.animation.Animation; android.graphics.Path; android.graphics.Paint; android.graphics.Canvas; android.graphics.DrawFilter; android.graphics.RectF'
$ws.Range("B5").Value = 'android.view.View; android.os.Build; android.support.v4.view.ViewPropertyAnimatorUpdateListener; android.util.Log; android.view.Gravity; android.graphics.Color; android.util.AttributeSet; android.content.res.TypedArray; android.view.MotionEvent; android.widget.AbsListView; android.support.v4.view.ViewCompat; android.widget.FrameLayout; android.support.v4.view.ViewPropertyAnimatorCompat; android.widget.Button; android.content.Context; android.view.ViewGroup; android.view.animation.DecelerateInterpolator'
$ws.Range("B6").Value = 'android.view.View; android.widget.AbsListView; android.os.Bundle; android.widget.ImageView; com.cjj.MaterialRefreshListener; android.view.ViewGroup; android.support.v7.widget.LinearLayoutManager; android.os.Handler; android.view.MenuItem; android.os.Build; android.support.v4.view.ViewPropertyAnimatorUpdateListener; android.graphics.Color; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; android.support.v4.view.ViewPropertyAnimatorCompat; android.support.v7.widget.Toolbar; java.util.Random; com.cjj.MaterialHeadListener; android.view.animation.DecelerateInterpolator; android.view.Menu; android.util.Log; android.support.v7.widget.DefaultItemAnimator; android.support.v7.widget.RecyclerView; android.widget.Toast; java.util.List; java.util.ArrayList; android.content.Context; android.view.Gravity; android.util.AttributeSet; android.content.res.TypedArray; android.view.MotionEvent; android.support.v4.view.ViewCompat; android.support.annotation.Nullable; android.widget.FrameLayout; android.widget.Button'
$ws.Range("B7").Value = 'android.os.Handler; com.cjj.MaterialRefreshLayout; android.os.Bundle; android.widget.ArrayAdapter; android.widget.Toast; android.widget.ListView; com.cjj.MaterialRefreshListener'
$ws.Range("B8").Value = 'android.view.View; android.widget.AbsListView; android.os.Bundle; android.widget.ImageView; android.widget.ListView; android.view.ViewGroup; com.cjj.MaterialRefreshListener; android.support.v7.widget.LinearLayoutManager; android.os.Handler; android.view.MenuItem; android.os.Build; android.support.v4.view.ViewPropertyAnimatorUpdateListener; android.graphics.Color; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; android.support.v4.view.ViewPropertyAnimatorCompat; android.support.v7.widget.Toolbar; java.util.Random; com.cjj.MaterialHeadListener; android.view.animation.DecelerateInterpolator; android.view.Menu; android.util.Log; android.support.v7.widget.DefaultItemAnimator; android.support.v7.widget.RecyclerView; android.widget.ArrayAdapter; android.widget.Toast; java.util.List; java.util.ArrayList; android.content.Context; android.view.Gravity; android.util.AttributeSet; android.content.res.TypedArray; android.view.MotionEvent; android.support.v4.view.ViewCompat; android.support.annotation.Nullable; android.widget.FrameLayout; android.widget.Button'
$ws.Range("B9").Value = 'android.view.View; android.graphics.RadialGradient; android.widget.AbsListView; android.animation.AnimatorSet; android.os.Bundle; android.widget.ImageView; android.widget.ListView; android.graphics.drawable.ShapeDrawable; android.view.ViewGroup; com.cjj.MaterialRefreshListener; android.support.v7.widget.LinearLayoutManager; android.view.MenuItem; android.os.Build; android.support.v4.view.ViewPropertyAnimatorUpdateListener; android.graphics.Color; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; android.support.v4.view.ViewPropertyAnimatorCompat; android.support.v7.widget.Toolbar; java.util.Random; com.cjj.MaterialHeadListener; android.view.animation.DecelerateInterpolator; android.view.Menu; android.util.Log; android.support.v7.widget.DefaultItemAnimator; android.content.res.Resources; android.support.v7.widget.RecyclerView; android.net.Uri; android.widget.ArrayAdapter; android.widget.Toast; java.util.List; java.util.ArrayList; android.content.Context; android.widget.Button; android.graphics.Shader; android.content.Intent; android.graphics.drawable.Drawable; android.graphics.drawable.shapes.OvalShape; android.view.Gravity; android.animation.ObjectAnimator; android.util.AttributeSet; android.content.res.TypedArray; android.view.MotionEvent; android.support.v4.view.ViewCompat; android.support.annotation.Nullable; android.view.animation.Animation; android.widget.FrameLayout; android.graphics.Paint; android.graphics.Canvas'
$ws.Range("B10").Value = 'android.view.View; android.view.animation.AccelerateDecelerateInterpolator; android.graphics.Paint.Style; android.widget.AbsListView; android.animation.AnimatorSet; android.graphics.PixelFormat; android.view.ViewGroup; android.graphics.Rect; android.os.Build; android.support.v4.view.ViewPropertyAnimatorUpdateListener; android.graphics.Color; android.support.v4.view.ViewPropertyAnimatorCompat; android.util.DisplayMetrics; android.view.animation.DecelerateInterpolator; android.view.animation.Transformation; android.graphics.ColorFilter; android.util.Log; android.view.animation.LinearInterpolator; android.content.res.Resources; android.view.animation.Interpolator; java.util.ArrayList; android.content.Context; android.graphics.drawable.Drawable; android.graphics.drawable.Animatable; android.view.Gravity; android.animation.ObjectAnimator; android.util.AttributeSet; android.content.res.TypedArray; android.view.MotionEvent; android.support.v4.view.ViewCompat; android.widget.FrameLayout; android.graphics.Path; android.view.animation.Animation; android.graphics.Paint; android.graphics.Canvas; android.graphics.RectF'
$ws.Range("B11").Value = 'android.view.View; android.widget.AbsListView; android.os.Bundle; android.widget.ImageView; android.widget.ListView; android.view.WindowManager; android.view.ViewGroup; com.cjj.MaterialRefreshListener; android.support.v7.widget.LinearLayoutManager; android.os.Build; android.view.MenuItem; android.graphics.Color; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; com.nineoldandroids.animation.Animator; android.support.v7.widget.Toolbar; com.nineoldandroids.view.ViewHelper; java.util.Random; com.cjj.MaterialHeadListener; android.view.animation.DecelerateInterpolator; android.view.Menu; android.util.Log; android.support.v7.widget.DefaultItemAnimator; android.support.v7.widget.RecyclerView; android.widget.ArrayAdapter; android.widget.Toast; android.view.animation.AnimationSet; java.util.ArrayList; android.content.Context; java.util.List; android.content.Intent; com.nineoldandroids.animation.AnimatorSet; android.view.Gravity; android.widget.LinearLayout; android.util.AttributeSet; android.content.res.TypedArray; com.nineoldandroids.animation.ObjectAnimator; android.view.MotionEvent; com.nineoldandroids.animation.ValueAnimator; android.support.v4.view.ViewCompat; android.view.animation.AccelerateInterpolator; android.widget.FrameLayout; android.support.annotation.Nullable; com.nineoldandroids.animation.AnimatorListenerAdapter'
$ws.Range("B12").Value = 'android.view.View; android.view.animation.AccelerateDecelerateInterpolator; android.graphics.RadialGradient; android.graphics.Paint.Style; android.widget.AbsListView; android.graphics.PixelFormat; android.os.Bundle; android.widget.ImageView; android.graphics.drawable.ShapeDrawable; android.view.ViewGroup; android.view.WindowManager; android.support.v7.widget.LinearLayoutManager; com.cjj.MaterialRefreshListener; android.graphics.Rect; android.os.Build; android.os.Handler; android.view.MenuItem; android.graphics.Color; android.app.Application; android.view.LayoutInflater; com.cjj.MaterialRefreshLayout; com.nineoldandroids.animation.Animator; android.view.animation.BounceInterpolator; android.support.v7.widget.Toolbar; android.util.DisplayMetrics; com.nineoldandroids.view.ViewHelper; java.lang.annotation.RetentionPolicy; java.util.Random; android.support.annotation.NonNull; android.view.animation.DecelerateInterpolator; android.test.ApplicationTestCase; android.view.animation.Transformation; android.graphics.ColorFilter; android.util.Log; android.view.animation.LinearInterpolator; android.support.v7.widget.DefaultItemAnimator; android.view.Menu; android.content.res.Resources; android.net.Uri; android.support.v7.widget.RecyclerView; android.view.animation.Interpolator; android.widget.Toast; android.view.animation.AnimationSet; java.util.ArrayList; android.content.Context; java.util.List; android.graphics.drawable.Drawable; android.graphics.Shader; com.nineoldandroids.animation.AnimatorSet; android.support.v7.app.AppCompatActivity; android.graphics.drawable.Animatable; android.graphics.drawable.shapes.OvalShape; java.lang.annotation.Retention; android.view.Gravity; android.widget.LinearLayout; android.util.AttributeSet; android.content.res.TypedArray; com.nineoldandroids.animation.ObjectAnimator; android.view.MotionEvent; com.nineoldandroids.animation.ValueAnimator; android.support.v4.view.ViewCompat; android.view.animation.Animation; android.view.animation.AccelerateInterpolator; android.widget.FrameLayout; android.graphics.Path; android.graphics.Paint; android.support.annotation.IntDef; android.graphics.Canvas; android.support.annotation.Nullable; com.nineoldandroids.animation.AnimatorListenerAdapter; android.graphics.RectF'
